$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("RBFF")

# Commit: "Take out biomass from building electrification"
# On the RBFF sheet, column G is the "biomass" source fuel and row headers
# down column A give the recipient fuel. Row 2 = electricity, row 7 = biomass.
# Previously 100% of biomass was redirected to electricity (G2=1, G7=0) as
# part of the building electrification fuel-shifting policy. Now biomass is
# excluded from that shift, so it maps back to itself instead (G2=0, G7=1).
$ws.Range("G2").Value = 0
$ws.Range("G7").Value = 1

# Restore the cell selection that was left on the About sheet.
$wsAbout.Activate()
$wsAbout.Range("J13").Select()

# Leave the RBFF sheet as the active/selected tab, with its prior selection.
$ws.Activate()
$ws.Range("M7").Select()
